$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 20000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -20924
$ws.Range("H98").Value = 3101.5652
$ws.Range("I98").Value = 2045.1538
$ws.Range("J98").Value = 8987.286
$ws.Range("K98").Value = 2045.1538
$ws.Range("L98").Value = 8987.286
$ws.Range("M98").Value = -547.1538
$ws.Range("N98").Value = -11983.286
$ws.Range("H121").Value = 1554.2858
$ws.Range("I121").Value = 970
$ws.Range("J121").Value = 2333.3333
$ws.Range("K121").Value = 2910
$ws.Range("L121").Value = 6999.999899999999
$ws.Range("M121").Value = -1163
$ws.Range("N121").Value = -10493.9999
$ws.Range("H122").Value = 3101.5652
$ws.Range("I122").Value = 2045.1538
$ws.Range("J122").Value = 8987.286
$ws.Range("K122").Value = 6135.4614
$ws.Range("L122").Value = 26961.858
$ws.Range("M122").Value = -3685.4614
$ws.Range("N122").Value = -31861.858
$ws.Range("H137").Value = 1148.8572
$ws.Range("I137").Value = 1040.9048
$ws.Range("J137").Value = 1472.7142
$ws.Range("K137").Value = 3122.7144
$ws.Range("L137").Value = 4418.142599999999
$ws.Range("M137").Value = -572.7143999999998
$ws.Range("N137").Value = -9518.142599999999
$ws.Range("H138").Value = 2454.5073
$ws.Range("I138").Value = 1447.1777
$ws.Range("J138").Value = 4343.25
$ws.Range("K138").Value = 4341.5331
$ws.Range("L138").Value = 13029.75
$ws.Range("M138").Value = 798.4669000000004
$ws.Range("N138").Value = -23309.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3054.4285
$ws.Range("I61").Value = 2154.05
$ws.Range("J61").Value = 5305.375
$ws.Range("K61").Value = 2154.05
$ws.Range("L61").Value = 5305.375
$ws.Range("M61").Value = -1942.05
$ws.Range("N61").Value = -5729.375
$ws.Range("H74").Value = 709.67566
$ws.Range("I74").Value = 595.82355
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 595.82355
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = 278.17645
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 709.67566
$ws.Range("I77").Value = 595.82355
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 2979.11775
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = 1388.88225
$ws.Range("N77").Value = -18736
$ws.Range("H122").Value = 1103.125
$ws.Range("I122").Value = 1258.7142
$ws.Range("J122").Value = 14
$ws.Range("K122").Value = 3776.1426
$ws.Range("L122").Value = 42
$ws.Range("M122").Value = -1326.1426
$ws.Range("N122").Value = -4942
$ws.Range("H132").Value = 8887.862999999999
$ws.Range("I132").Value = 17941.875
$ws.Range("J132").Value = 3714.1428
$ws.Range("K132").Value = 53825.625
$ws.Range("L132").Value = 11142.4284
$ws.Range("M132").Value = -51295.625
$ws.Range("N132").Value = -16202.4284
$ws.Range("H136").Value = 3054.4285
$ws.Range("I136").Value = 2154.05
$ws.Range("J136").Value = 5305.375
$ws.Range("K136").Value = 6462.150000000001
$ws.Range("L136").Value = 15916.125
$ws.Range("M136").Value = -3912.150000000001
$ws.Range("N136").Value = -21016.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 157107.84
$ws.Range("I86").Value = 4085
$ws.Range("J86").Value = 335634.5
$ws.Range("K86").Value = 4085
$ws.Range("L86").Value = 335634.5
$ws.Range("M86").Value = -2962
$ws.Range("N86").Value = -337880.5
$ws.Range("H89").Value = 157107.84
$ws.Range("I89").Value = 4085
$ws.Range("J89").Value = 335634.5
$ws.Range("K89").Value = 20425
$ws.Range("L89").Value = 1678172.5
$ws.Range("M89").Value = -14809
$ws.Range("N89").Value = -1689404.5
$ws.Range("H107").Value = 86401.664
$ws.Range("I107").Value = 103282
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 103282
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -101362
$ws.Range("N107").Value = -5840
$ws.Range("H134").Value = 1631.5077
$ws.Range("I134").Value = 1505.661
$ws.Range("J134").Value = 2869
$ws.Range("K134").Value = 4516.983
$ws.Range("L134").Value = 8607
$ws.Range("M134").Value = -1981.983
$ws.Range("N134").Value = -13677

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1430.2712
$ws.Range("I31").Value = 914
$ws.Range("J31").Value = 3960
$ws.Range("K31").Value = 914
$ws.Range("L31").Value = 3960
$ws.Range("M31").Value = -619
$ws.Range("N31").Value = -4550
$ws.Range("H34").Value = 1430.2712
$ws.Range("I34").Value = 914
$ws.Range("J34").Value = 3960
$ws.Range("K34").Value = 914
$ws.Range("L34").Value = 3960
$ws.Range("M34").Value = -712
$ws.Range("N34").Value = -4364
$ws.Range("H58").Value = 1364.1316
$ws.Range("I58").Value = 1342.3334
$ws.Range("J58").Value = 1417.6364
$ws.Range("K58").Value = 1342.3334
$ws.Range("L58").Value = 1417.6364
$ws.Range("M58").Value = -1139.3334
$ws.Range("N58").Value = -1823.6364
$ws.Range("H132").Value = 1614.8462
$ws.Range("I132").Value = 1236.3448
$ws.Range("J132").Value = 2712.5
$ws.Range("K132").Value = 3709.0344
$ws.Range("L132").Value = 8137.5
$ws.Range("M132").Value = -1179.0344
$ws.Range("N132").Value = -13197.5
$ws.Range("H134").Value = 1908.625
$ws.Range("I134").Value = 1200.6207
$ws.Range("J134").Value = 3775.182
$ws.Range("K134").Value = 3601.8621
$ws.Range("L134").Value = 11325.546
$ws.Range("M134").Value = -1066.8621
$ws.Range("N134").Value = -16395.546
$ws.Range("H136").Value = 1364.1316
$ws.Range("I136").Value = 1342.3334
$ws.Range("J136").Value = 1417.6364
$ws.Range("K136").Value = 4027.0002
$ws.Range("L136").Value = 4252.9092
$ws.Range("M136").Value = -1477.0002
$ws.Range("N136").Value = -9352.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1263.1666
$ws.Range("I122").Value = 930
$ws.Range("J122").Value = 1596.3334
$ws.Range("K122").Value = 8370
$ws.Range("L122").Value = 14367.0006
$ws.Range("M122").Value = -5920
$ws.Range("N122").Value = -19267.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 2797.7896
$ws.Range("I122").Value = 2240
$ws.Range("J122").Value = 3417.5557
$ws.Range("K122").Value = 6720
$ws.Range("L122").Value = 10252.6671
$ws.Range("M122").Value = -4270
$ws.Range("N122").Value = -15152.6671
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 2746.4285
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 2986.3635
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 8959.0905
$ws.Range("M126").Value = -3129.9998
$ws.Range("N126").Value = -13899.0905
$ws.Range("H132").Value = 2563.2156
$ws.Range("I132").Value = 2364.2273
$ws.Range("J132").Value = 3814
$ws.Range("K132").Value = 7092.6819
$ws.Range("L132").Value = 11442
$ws.Range("M132").Value = -4562.6819
$ws.Range("N132").Value = -16502

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1768.5264
$ws.Range("I16").Value = 1682.3529
$ws.Range("K16").Value = 1682.3529
$ws.Range("M16").Value = -1512.3529
$ws.Range("H69").Value = 500163
$ws.Range("J69").Value = 500163
$ws.Range("L69").Value = 500163
$ws.Range("N69").Value = -501785
$ws.Range("H72").Value = 500163
$ws.Range("J72").Value = 500163
$ws.Range("L72").Value = 1500489
$ws.Range("N72").Value = -1508601
$ws.Range("H122").Value = 16670209
$ws.Range("I122").Value = 12503326
$ws.Range("J122").Value = 28575586
$ws.Range("K122").Value = 37509978
$ws.Range("L122").Value = 85726758
$ws.Range("M122").Value = -37507528
$ws.Range("N122").Value = -85731658
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H136").Value = 16292022
$ws.Range("I136").Value = 20001328
$ws.Range("K136").Value = 60003984
$ws.Range("M136").Value = -60001434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23114
$ws.Range("J64").Value = 23114
$ws.Range("L64").Value = 23114
$ws.Range("N64").Value = -23610
$ws.Range("H67").Value = 23114
$ws.Range("J67").Value = 23114
$ws.Range("L67").Value = 23114
$ws.Range("N67").Value = -24830
$ws.Range("H132").Value = 2393.3572
$ws.Range("I132").Value = 1959.0834
$ws.Range("K132").Value = 5877.2502
$ws.Range("M132").Value = -3347.2502
$ws.Range("H136").Value = 1358.7
$ws.Range("I136").Value = 1392.4828
$ws.Range("J136").Value = 379
$ws.Range("K136").Value = 4177.4484
$ws.Range("L136").Value = 1137
$ws.Range("M136").Value = -1627.4484
$ws.Range("N136").Value = -6237
